$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 45.71598933333333
$ws.Range("H2").Value = 137.147968
$ws.Range("I2").Value = 0.6549002937372808
$ws.Range("J2").Value = 0.6549002937372808
$ws.Range("M2").Value = 449.104309
$ws.Range("N2").Value = 1347.312927
$ws.Range("O2").Value = 0.9710020245482639
$ws.Range("P2").Value = 0.9710020245482639
$ws.Range("Q2").Value = 20531.24779979804
$ws.Range("R2").Value = 184781.2301981823
$ws.Range("S2").Value = 0.6359095110961523
$ws.Range("T2").Value = 0.6359095110961523
# Row 3
$ws.Range("G3").Value = 45.71598933333333
$ws.Range("H3").Value = 137.147968
$ws.Range("I3").Value = 0.6549002937372808
$ws.Range("J3").Value = 0.6549002937372808
$ws.Range("O3").Value = 0.01131353526791385
$ws.Range("P3").Value = 0.01131353526791385
$ws.Range("Q3").Value = 239.2178288045866
$ws.Range("R3").Value = 2152.96045924128
$ws.Range("S3").Value = 0.007409237570163862
$ws.Range("T3").Value = 0.007409237570163864
# Row 4
$ws.Range("G4").Value = 45.71598933333333
$ws.Range("H4").Value = 137.147968
$ws.Range("I4").Value = 0.6549002937372808
$ws.Range("J4").Value = 0.6549002937372808
$ws.Range("M4").Value = 5.266527
$ws.Range("N4").Value = 15.799581
$ws.Range("O4").Value = 0.0113866829528418
$ws.Range("P4").Value = 0.0113866829528418
$ws.Range("Q4").Value = 240.764492155712
$ws.Range("R4").Value = 2166.880429401408
$ws.Range("S4").Value = 0.007457142010509381
$ws.Range("T4").Value = 0.007457142010509381
# Row 5
$ws.Range("G5").Value = 45.71598933333333
$ws.Range("H5").Value = 137.147968
$ws.Range("I5").Value = 0.6549002937372808
$ws.Range("J5").Value = 0.6549002937372808
$ws.Range("M5").Value = 2.912815666666667
$ws.Range("N5").Value = 8.738447000000001
$ws.Range("O5").Value = 0.006297757230980464
$ws.Range("P5").Value = 0.006297757230980464
$ws.Range("Q5").Value = 133.1622499472996
$ws.Range("R5").Value = 1198.460249525696
$ws.Range("S5").Value = 0.00412440306045519
$ws.Range("T5").Value = 0.00412440306045519
# Row 6
$ws.Range("I6").Value = 0.1818108415648851
$ws.Range("J6").Value = 0.1818108415648851
$ws.Range("M6").Value = 449.104309
$ws.Range("N6").Value = 1347.312927
$ws.Range("O6").Value = 0.9710020245482639
$ws.Range("P6").Value = 0.9710020245482639
$ws.Range("Q6").Value = 5699.804193943337
$ws.Range("R6").Value = 51298.23774549003
$ws.Range("S6").Value = 0.1765386952443271
$ws.Range("T6").Value = 0.176538695244327
# Row 7
$ws.Range("I7").Value = 0.1818108415648851
$ws.Range("J7").Value = 0.1818108415648851
$ws.Range("O7").Value = 0.01131353526791385
$ws.Range("P7").Value = 0.01131353526791385
$ws.Range("S7").Value = 0.002056923368133424
$ws.Range("T7").Value = 0.002056923368133424
# Row 8
$ws.Range("I8").Value = 0.1818108415648851
$ws.Range("J8").Value = 0.1818108415648851
$ws.Range("M8").Value = 5.266527
$ws.Range("N8").Value = 15.799581
$ws.Range("O8").Value = 0.0113866829528418
$ws.Range("P8").Value = 0.0113866829528418
$ws.Range("Q8").Value = 66.840090554811
$ws.Range("R8").Value = 601.5608149932989
$ws.Range("S8").Value = 0.002070222410288698
$ws.Range("T8").Value = 0.002070222410288698
# Row 9
$ws.Range("I9").Value = 0.1818108415648851
$ws.Range("J9").Value = 0.1818108415648851
$ws.Range("M9").Value = 2.912815666666667
$ws.Range("N9").Value = 8.738447000000001
$ws.Range("O9").Value = 0.006297757230980464
$ws.Range("P9").Value = 0.006297757230980464
$ws.Range("Q9").Value = 36.96797964379034
$ws.Range("R9").Value = 332.711816794113
$ws.Range("S9").Value = 0.001145000542135899
$ws.Range("T9").Value = 0.001145000542135898
# Row 10
$ws.Range("G10").Value = 11.24784666666667
$ws.Range("H10").Value = 33.74354
$ws.Range("I10").Value = 0.161130015850732
$ws.Range("J10").Value = 0.161130015850732
$ws.Range("M10").Value = 449.104309
$ws.Range("N10").Value = 1347.312927
$ws.Range("O10").Value = 0.9710020245482639
$ws.Range("P10").Value = 0.9710020245482639
$ws.Range("Q10").Value = 5051.456404971286
$ws.Range("R10").Value = 45463.10764474157
$ws.Range("S10").Value = 0.1564575716065546
$ws.Range("T10").Value = 0.1564575716065546
# Row 11
$ws.Range("G11").Value = 11.24784666666667
$ws.Range("H11").Value = 33.74354
$ws.Range("I11").Value = 0.161130015850732
$ws.Range("J11").Value = 0.161130015850732
$ws.Range("O11").Value = 0.01131353526791385
$ws.Range("P11").Value = 0.01131353526791385
$ws.Range("Q11").Value = 58.85655101343333
$ws.Range("R11").Value = 529.7089591208999
$ws.Range("S11").Value = 0.001822950117046773
$ws.Range("T11").Value = 0.001822950117046773
# Row 12
$ws.Range("G12").Value = 11.24784666666667
$ws.Range("H12").Value = 33.74354
$ws.Range("I12").Value = 0.161130015850732
$ws.Range("J12").Value = 0.161130015850732
$ws.Range("M12").Value = 5.266527
$ws.Range("N12").Value = 15.799581
$ws.Range("O12").Value = 0.0113866829528418
$ws.Range("P12").Value = 0.0113866829528418
$ws.Range("Q12").Value = 59.23708816185999
$ws.Range("R12").Value = 533.1337934567399
$ws.Range("S12").Value = 0.001834736404678658
$ws.Range("T12").Value = 0.001834736404678658
# Row 13
$ws.Range("G13").Value = 11.24784666666667
$ws.Range("H13").Value = 33.74354
$ws.Range("I13").Value = 0.161130015850732
$ws.Range("J13").Value = 0.161130015850732
$ws.Range("M13").Value = 2.912815666666667
$ws.Range("N13").Value = 8.738447000000001
$ws.Range("O13").Value = 0.006297757230980464
$ws.Range("P13").Value = 0.006297757230980464
$ws.Range("Q13").Value = 32.76290398693111
$ws.Range("R13").Value = 294.86613588238
$ws.Range("S13").Value = 0.001014757722451944
$ws.Range("T13").Value = 0.001014757722451944
# Row 14
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.1507006666666667
$ws.Range("H14").Value = 0.452102
$ws.Range("I14").Value = 0.00215884884710222
$ws.Range("J14").Value = 0.00215884884710222
$ws.Range("M14").Value = 449.104309
$ws.Range("N14").Value = 1347.312927
$ws.Range("O14").Value = 0.9710020245482639
$ws.Range("P14").Value = 0.9710020245482639
$ws.Range("Q14").Value = 67.68031876917267
$ws.Range("R14").Value = 609.1228689225539
$ws.Range("S14").Value = 0.002096246601229941
$ws.Range("T14").Value = 0.002096246601229941
# Row 15
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.1507006666666667
$ws.Range("H15").Value = 0.452102
$ws.Range("I15").Value = 0.00215884884710222
$ws.Range("J15").Value = 0.00215884884710222
$ws.Range("O15").Value = 0.01131353526791385
$ws.Range("P15").Value = 0.01131353526791385
$ws.Range("Q15").Value = 0.7885706249633333
$ws.Range("R15").Value = 7.09713562467
$ws.Range("S15").Value = [double]"2.442421256978611E-05"
$ws.Range("T15").Value = [double]"2.442421256978611E-05"
# Row 16
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.1507006666666667
$ws.Range("H16").Value = 0.452102
$ws.Range("I16").Value = 0.00215884884710222
$ws.Range("J16").Value = 0.00215884884710222
$ws.Range("M16").Value = 5.266527
$ws.Range("N16").Value = 15.799581
$ws.Range("O16").Value = 0.0113866829528418
$ws.Range("P16").Value = 0.0113866829528418
$ws.Range("Q16").Value = 0.7936691299180001
$ws.Range("R16").Value = 7.143022169262
$ws.Range("S16").Value = [double]"2.458212736506102E-05"
$ws.Range("T16").Value = [double]"2.458212736506101E-05"
# Row 17
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.1507006666666667
$ws.Range("H17").Value = 0.452102
$ws.Range("I17").Value = 0.00215884884710222
$ws.Range("J17").Value = 0.00215884884710222
$ws.Range("M17").Value = 2.912815666666667
$ws.Range("N17").Value = 8.738447000000001
$ws.Range("O17").Value = 0.006297757230980464
$ws.Range("P17").Value = 0.006297757230980464
$ws.Range("Q17").Value = 0.4389632628437778
$ws.Range("R17").Value = 3.950669365594
$ws.Range("S17").Value = [double]"1.359590593743184E-05"
$ws.Range("T17").Value = [double]"1.359590593743184E-05"
